# Daily limits fix: shrink the "...-7" / "...,7" weekday ranges in column C
# down to "...-6" / "...,6" for the rows listed below, and highlight the
# affected rows with a yellow fill so the change is visually obvious.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column C
$updates = @{
    122 = "1-6"
    129 = "1-6"
    133 = "1-6"
    135 = "1-6"
    136 = "1-6"
    137 = "1-6"
    138 = "3-6"
    140 = "2-6"
    141 = "2-6"
    142 = "2-6"
    143 = "1,6"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
    $ws.Range("A" + $row + ":E" + $row).Interior.Color = 65535
}

# Restore the view/selection to where the edited rows are, matching the
# saved state of the workbook after the edit.
$ws.Range("C143").Select()
